$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Fix the output creator loop: columns H:K for rows 15-26 should all be 20
# (previously 15 for rows 15-20 and 8 for rows 21-26).
$ws.Range("H15:K26").Value = 20

# Update the active selection on the Data sheet to L7
$ws.Activate()
$ws.Range("L7").Select()
